$d = $word.ActiveDocument

# Locate the paragraph that holds "LOQ4031: Química Geral I (Requisito fraco)".
$loq = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r") -eq "LOQ4031: Química Geral I (Requisito fraco)") {
        $loq = $p
        break
    }
}

if ($loq -ne $null) {
    # Right after it come three paragraphs that must be removed entirely:
    #   1) an empty "Normal" paragraph
    #   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
    #   3) "© 2020 . Contact: ... Creative Commons Attribution"
    # The paragraph that originally followed those (also an empty "Normal"
    # paragraph, right before the page-break paragraph) must stay untouched.
    $p1 = $loq.Next()
    $p2 = $p1.Next()
    $p3 = $p2.Next()

    $start = $p1.Range.Start
    $finish = $p3.Range.End

    $d.Range($start, $finish).Delete()
}
